# 7.1.2.xlsx update:
#  - Rewrite the indicator description / unit labels (B3, B4)
#  - Shrink row 4's height now that "Prozent" is a one-line label
#  - Insert a new "Anmerkung:" / "Alle Daten geschätzt." row right above the
#    "Copyright:" row (pushes Copyright/Vervielfältigung/blank rows down by one)
#  - Nudge the small footnote icon picture down to stay next to the
#    "Copyright:" row after the insert
#  - Update the selected cell to match the saved view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header / unit text updates -------------------------------------------------
$ws.Range("B3").Value = "Anteil der Bevölkerung, der vorwiegend saubere Energieträger und Technologien nutzt"
$ws.Range("B4").Value = "Prozent"
$ws.Rows(4).RowHeight = 15

# --- Insert the new "Anmerkung:" row ---------------------------------------------
$ws.Rows("37:37").Insert()
$ws.Range("A37").Value = "Anmerkung:"
$ws.Range("B37").Value = "Alle Daten geschätzt."

# --- Keep the small footnote picture aligned with the (now shifted) Copyright row
$icon = $ws.Shapes.Item("Grafik 6")
$icon.Top = $icon.Top() + $ws.Rows(37).RowHeight()

# --- Selection bookkeeping ---------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("I39").Select()
